$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$headers = @("Código", "Nome", "Descrição", "Categoria", "Unidade", "Marca", "Grupo")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $cell.Value = $headers[$i]
    $cell.Font.Bold = $true
    $cell.Font.Size = 12
}

# Data rows
$data = @(
    @("s35", "SORVETE NATA", "SORVETE NATA (4X1,5L)", "LPC- Sorv PT 1,5L NOBRE BCD", "Caixa", "Bariloche", "Sorvete"),
    @("S2626", "SORVETE 1,5L MOUSSE LIMAO", "SORVETE 1,5L MOUSSE LIMAO (4X1,5L)", "LPC- Sorv PT 1,8L MESCLADO BCD", "Unidade", "Bariloche", "Sorvete"),
    @("S262", "SORVETE FRUTAS VERM", "SORVETE FRUTAS VERM (4X1,5L)", "LPC- Sorv PT 1,8L MESCLADO BCD", "Caixa", "Bariloche", "Sorvete")
)

for ($r = 0; $r -lt $data.Length; $r++) {
    $row = $data[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $ws.Cells.Item($r + 2, $c + 1).Value = $row[$c]
    }
}

for ($c = 1; $c -le 7; $c++) {
    $ws.Columns.Item($c).AutoFit()
}
$ws.Range("G1").Select() | Out-Null
